$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 262, shifting existing rows 262:346 down to 263:347
$ws.Rows.Item(262).Insert()

# Populate the new row 262 with the new data record
$ws.Range("A262").Value = 4
$ws.Range("B262").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C262").Value = "Los Lagos"
$ws.Range("D262").Value = 44627
$ws.Range("E262").Value = 10
$ws.Range("F262").Value = 100114001
$ws.Range("G262").Value = "Papa"
$ws.Range("H262").Value = "Patagonia"
$ws.Range("I262").Value = "1a (cosecha)"
$ws.Range("J262").Value = 300
$ws.Range("K262").Value = 7000
$ws.Range("L262").Value = 7000
$ws.Range("M262").Value = 7000
$ws.Range("N262").Value = '$/saco 25 kilos'
$ws.Range("O262").Value = "Provincia de Llanquihue"
$ws.Range("P262").Value = 280
$ws.Range("Q262").Value = 25
$ws.Range("R262").Value = "Hortaliza"

# Match the date-number style used by the rest of column D
$ws.Range("D262").NumberFormat = $ws.Range("D263").NumberFormat
